# [30-10-2023][Rey] Final Progress Pre Paper
# Refresh the AlgoEdge results table: update the "edge" counts (col C) and
# re-index the per-algorithm row numbering (col A) for the Spinglass,
# Girvan Newman and Belief groups, and drop the two trailing Belief rows
# that no longer exist in the refreshed run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two now-unused trailing rows (old rows 29-30), shrinking the
# sheet from A1:C30 down to A1:C28.
$ws.Rows("29:30").Delete()

# --- Spinglass group (rows 2-15, index column already 0-13) -------------
$ws.Range("C3").Value = 152
$ws.Range("C4").Value = 151
$ws.Range("C5").Value = 116
$ws.Range("C6").Value = 151
$ws.Range("C7").Value = 102
$ws.Range("C8").Value = 99
$ws.Range("C9").Value = 61
$ws.Range("C11").Value = 79
$ws.Range("C12").Value = 80
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 2

# --- Girvan Newman group (rows 16-21) now starts at row 16 instead of 18,
# and has two extra leading rows (index 0 and 1) that used to belong to
# the tail of the Spinglass block.
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = "Girvan Newman"
$ws.Range("C16").Value = 420

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Girvan Newman"
$ws.Range("C17").Value = 266

$ws.Range("A18").Value = 2
$ws.Range("C18").Value = 273

$ws.Range("A19").Value = 3
$ws.Range("C19").Value = 192

$ws.Range("A20").Value = 4
$ws.Range("C20").Value = 128

$ws.Range("A21").Value = 5
$ws.Range("C21").Value = 108

# --- Belief group (rows 22-28) now starts at row 22 instead of 24, and
# has two extra leading rows (index 0 and 1) pulled in from the old
# Girvan Newman tail.
$ws.Range("A22").Value = 0
$ws.Range("B22").Value = "Belief"
$ws.Range("C22").Value = 230

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Belief"
$ws.Range("C23").Value = 232

$ws.Range("A24").Value = 2
$ws.Range("C24").Value = 220

$ws.Range("A25").Value = 3
$ws.Range("C25").Value = 222

$ws.Range("A26").Value = 4
$ws.Range("C26").Value = 234

$ws.Range("A27").Value = 5
$ws.Range("C27").Value = 105

$ws.Range("A28").Value = 6
$ws.Range("C28").Value = 116
